$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-28 04:06:32"
$wsZhCn.Range("G3").Value = "2016-01-28 04:07:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-28 04:06:42"
$wsDeDe.Range("G3").Value = "2016-01-28 04:07:31"
